# Add four new "planned" Portfolio Optimizer & Constructor sprint tasks
# (S16_G03_TB00X, S16_G03_TB00Y, S16_G03_TF00X, S16_G03_TB00Z) as new rows
# 156-159 on the sprint tracker sheet, matching the existing plain-default
# formatting used by the unstyled rows already in the sheet (e.g. 153-155).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$group = "Portfolio Optimizer & Constructor"

$rows = @(
    @{ Row = 156; Sprint = "S16"; Group = "G03"; GroupDesc = $group; Task = "S16_G03_TB00X";
       Desc = "Add configurable risk-free rate for portfolios and use it for Sharpe/Max-Sharpe calculations (construction + backtests).";
       Remarks = "Planned: risk_free_rate stored in portfolio risk profile and surfaced on Construction tab.";
       Status = "planned" },
    @{ Row = 157; Sprint = "S16"; Group = "G03"; GroupDesc = $group; Task = "S16_G03_TB00Y";
       Desc = "Add optimisation risk-profile presets (Conservative/Moderate/Aggressive) that map to sensible constraint defaults for Indian equities.";
       Remarks = "Planned: presets populate min/max weight, target volatility, max beta, and turnover limit.";
       Status = "planned" },
    @{ Row = 158; Sprint = "S16"; Group = "G03"; GroupDesc = $group; Task = "S16_G03_TF00X";
       Desc = "Show a small equity-curve preview for the optimised portfolio (run a short portfolio backtest behind the Construction tab).";
       Remarks = "Planned: reuse PortfolioService to simulate equity using current constraints and rebalance policy.";
       Status = "planned" },
    @{ Row = 159; Sprint = "S16"; Group = "G03"; GroupDesc = $group; Task = "S16_G03_TB00Z";
       Desc = "Add advanced factor/risk lookback presets (Short/Medium/Long) and plumb them into factor and risk model services.";
       Remarks = "Planned: expose as advanced options while keeping default as current 180–252 day lookbacks.";
       Status = "planned" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $values = @($r.Sprint, $r.Group, $r.GroupDesc, $r.Task, $r.Desc, $r.Remarks, $r.Status)
    for ($col = 1; $col -le 7; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        # Reset alignment to the sheet's plain default (general/bottom, no
        # wrap) instead of inheriting the bold/wrap column style, so these
        # new cells pick up the same "no explicit style" formatting as the
        # sheet's other unstyled data rows.
        $cell.WrapText = $false
        $cell.VerticalAlignment = -4107
        $cell.HorizontalAlignment = 1
        $cell.Value = $values[$col - 1]
    }
}
